$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new dSF (column F) value, per diff
$updates = @{
    2  = -4
    3  = -2
    4  = -2
    5  = -1
    6  = 5
    7  = 2
    8  = -1
    9  = -3
    10 = 3
    11 = 1
    12 = -1
    13 = -3
    14 = -1
    15 = 4
    16 = 2
    17 = 1
    18 = 3
    19 = 3
    20 = 1
    21 = -1
    22 = 3
    24 = -3
    25 = 2
    26 = 1
    27 = -6
    28 = -2
    29 = 3
    30 = 4
    31 = 2
    32 = -2
    33 = 4
    34 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
